$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a "feature table": columns A-F are
# set / feature / presenting_name / left_end / right_end / question(_heb).
# This edit:
#   1) renames the F1 header "question" -> "question_heb"
#   2) adds a new column G "question_heb_max_min_ideal" with per-row
#      Hebrew question templates (parameterised with "{}") for the four
#      "set A" rows, and simple duplicates of column F for the "set B/C" rows

# --- 1. Header row ---
$ws.Range("F1").Value = "question_heb"

# --- 2. New column G content for rows 2-5 (the "A" / set rows) ---
$ws.Range("G2").Value = "מה ההעדפה הפוליטית של {}"
$ws.Range("G3").Value = "איפה מתרחשים התחביבים של {}"
$ws.Range("G4").Value = "מה מבנה הגוף של {}"
$ws.Range("G5").Value = "מה רמת האינטליגנציה של {}"

# New column header - set after the row2-5 values and bolded like the
# rest of the header row.
$ws.Range("G1").Value = "question_heb_max_min_ideal"
$ws.Range("G1").Font.Bold = $true

# --- 3. Rows 6-10 (the "B"/"C" set rows): column G mirrors column F ---
$ws.Range("G6").Value = $ws.Range("F6").Value2
$ws.Range("G7").Value = $ws.Range("F7").Value2
$ws.Range("G8").Value = $ws.Range("F8").Value2
$ws.Range("G9").Value = $ws.Range("F9").Value2
$ws.Range("G10").Value = $ws.Range("F10").Value2

# --- 4. New column G should be as wide as column F ---
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# --- 5. Move/leave the active selection on the new header cell ---
$ws.Range("G1").Select() | Out-Null
